$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Remove the "bean factory" column (old column G) from the mapping table on Sheet1.
# This shifts the "oneWay"/"convertMethodAB" columns left from H/I to G/H and
# drops the now-unused shared strings ("classABeanFactory",
# "Class A \nbean factory\nclass", "org.dozer.factory.XMLBeanFactory").
$ws1.Columns.Item(7).Delete()

# The extra-tall header row (needed for the wrapped "bean factory" header) is no
# longer needed now that column is gone - let it size back to the default.
$ws1.Rows.Item(18).AutoFit()

# Sheet2 becomes the active/selected sheet (tab) instead of Sheet1.
$ws2.Activate()

# Update the remembered selections on each sheet.
$ws1.Range("G20").Select()
$ws2.Range("E11").Select()

$wb.Save()
